$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells to Text format so values like "14.60" or "1.00"
# keep their exact display (matches the source workbook, which stores these as text).
$priceCells = "D2","D3","D5","D6","D7","D10","D11","D12","D14","D15","D16","D17","D19","D20","D22","D23","D24","D25","D28","D30","D32","D34","D35","D37","D42","D43","D44","D45","D48","D49","D51"
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '42.509.17'
$ws.Range("E2").Value = '  +2.01%  '
$ws.Range("D3").Value = '2.287.25'
$ws.Range("E3").Value = '  +1.32%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '157.72'
$ws.Range("E5").Value = '  +15,645.90%  '
$ws.Range("D6").Value = '306.60'
$ws.Range("E6").Value = '  +1.08%  '
$ws.Range("D7").Value = '96.54'
$ws.Range("E7").Value = '  +5.71%  '
$ws.Range("E8").Value = '  +0.61%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").Value = '0.496'
$ws.Range("E10").Value = '  +3.69%  '
$ws.Range("D11").Value = '36.44'
$ws.Range("E11").Value = '  +13.58%  '
$ws.Range("D12").Value = '0.0804'
$ws.Range("E12").Value = '  +0.99%  '
$ws.Range("E13").Value = '  -1.49%  '
$ws.Range("D14").Value = '6.73'
$ws.Range("E14").Value = '  +2.59%  '
$ws.Range("D15").Value = '2.641.18'
$ws.Range("E15").Value = '  +1.31%  '
$ws.Range("D16").Value = '14.60'
$ws.Range("E16").Value = '  +2.99%  '
$ws.Range("D17").Value = '2.266.61'
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("E18").Value = '  +6.06%  '
$ws.Range("D19").Value = '42.406.05'
$ws.Range("E19").Value = '  +1.95%  '
$ws.Range("D20").Value = '12.83'
$ws.Range("E20").Value = '  +4.42%  '
$ws.Range("E21").Value = '  +1.97%  '
$ws.Range("D22").Value = '6.00'
$ws.Range("E22").Value = '  +2.02%  '
$ws.Range("D23").Value = '67.92'
$ws.Range("E23").Value = '  +2.05%  '
$ws.Range("D24").Value = '243.16'
$ws.Range("E24").Value = '  +1.29%  '
$ws.Range("D25").Value = '2.60'
$ws.Range("E25").Value = '  +1.37%  '
$ws.Range("E26").Value = '  +2.82%  '
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("D28").Value = '24.03'
$ws.Range("E28").Value = '  +0.47%  '
$ws.Range("E29").Value = '  +7.31%  '
$ws.Range("D30").Value = '9.61'
$ws.Range("E30").Value = '  +1.48%  '
$ws.Range("E31").Value = '  +2.30%  '
$ws.Range("D32").Value = '161.78'
$ws.Range("E32").Value = '  +1.05%  '
$ws.Range("E33").Value = '  +3.91%  '
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").Value = '0.0753'
$ws.Range("E35").Value = '  +1.72%  '
$ws.Range("E36").Value = '  +3.14%  '
$ws.Range("D37").Value = '17.42'
$ws.Range("E37").Value = '  +5.07%  '
$ws.Range("E38").Value = '  +4.38%  '
$ws.Range("E39").Value = '  +5.46%  '
$ws.Range("E40").Value = '  +0.36%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").Value = '4.19'
$ws.Range("E42").Value = '  +7.35%  '
$ws.Range("D43").Value = '2.35'
$ws.Range("E43").Value = '  +15.54%  '
$ws.Range("D44").Value = '2.006.09'
$ws.Range("E44").Value = '  -2.28%  '
$ws.Range("D45").Value = '19.31'
$ws.Range("E45").Value = '  -1.06%  '
$ws.Range("E46").Value = '  +2.72%  '
$ws.Range("E47").Value = '  +5.90%  '
$ws.Range("D48").Value = '10.20'
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("D49").Value = '54.06'
$ws.Range("E49").Value = '  +4.97%  '
$ws.Range("E50").Value = '  +1.93%  '
$ws.Range("D51").Value = '72.79'
$ws.Range("E51").Value = '  +0.47%  '
